$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "('Aswan Jaguar', ['{1}{G}{G}', 'Summon Jaguar', 'When Aswan Jaguar comes into play, choose a random creature type from those in target opponent" + [char]0x2019 + "s deck.', '{G}{G}, {T}: Bury target creature of the chosen type.', '2/2'])"

$ws.Range("A2").Value = $newValue

$ws.Range("A3:A7").EntireRow.Delete()
